# Apply the Spanish resume text edits described by the commit diff.
# Each call uses Find/Replace (Word COM interop) scoped to the whole
# document body so the exact phrase is matched and swapped in place.
#
# Find.Execute signature used below:
#   FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#   MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace

$d = $word.ActiveDocument

$edits = @(
    @{ Old = "Currículum vítae - Alex Wilber"; New = "Currículum vítae: Alex Wilber" },
    @{ Old = "Animación de Spark: Diseñador de animaciones (enero de 2021 - Presente)"; New = "Spark Animation: Diseñador de animaciones (enero de 2021 - actualidad)" },
    @{ Old = "Líder de un equipo de 12 animadores para crear animaciones 3D de alta calidad para diversos proyectos, como largometrajes, anuncios y videojuegos."; New = "Dirige un equipo de 12 animadores para crear animaciones 3D de alta calidad para diversos proyectos, como largometrajes, anuncios y videojuegos." },
    @{ Old = "Pixel Studio: Diseñador de animaciones (junio de 2018 - dic 2020)"; New = "Pixel Studio: Diseñador de animaciones (junio de 2018 - diciembre de 2020)" },
    @{ Old = "Animación flash: Diseñador de animaciones junior (septiembre de 2016 - mayo de 2018)"; New = "Flash Animation: Diseñador de animaciones junior (septiembre de 2016 - mayo de 2018)" },
    @{ Old = "Education"; New = "Educación" },
    @{ Old = "Maestro de Artes en Animación, Graduación esperada: dic 2025"; New = "Maestría en Humanidades con especialización en animación, graduación esperada: diciembre de 2025" },
    @{ Old = "El arte de la animación 3D: una guía para principiantes."; New = "The Art of 3D Animation: A Guide for Beginners." }
)

foreach ($edit in $edits) {
    $found = $d.Content.Find.Execute($edit.Old, $true, $false, $false, $false, $false,
                                      $true, 1, $false, $edit.New, 2)
    if (-not $found) {
        Write-Output "WARNING: phrase not found -> $($edit.Old)"
    }
}
